# Daily attendance processing - 2026-01-18 18:41:37
#
# In the "Recorded By" column (G), swap the display order of the two
# recorder names wherever they appear together, e.g.
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
#
# Cells in column G that contain only "System" or only
# "dnasr281@gmail.com" (no comma-joined pair) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Restrict the operation to column G ("Recorded By") and do a whole-cell
# (not partial) text replace so only exact matches are touched.
$col = $ws.Columns.Item(7)
$col.Replace($oldValue, $newValue, 1, 1, $false)
